$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Update "Bad Drivers" table data rows (3,4,5) ---
$ws.Range("C3").Value2 = 618
$ws.Range("D3").Value2 = 86.4
$ws.Range("C4").Value2 = 305
$ws.Range("D4").Value2 = 95.5
$ws.Range("B5").Value2 = 5
$ws.Range("C5").Value2 = 553
$ws.Range("D5").Value2 = 98

# --- Step 2: Remove the obsolete driver row (old row 6: AX201 23.10.0.8) ---
# This shifts everything below up by one row, which also aligns the
# "Good Drivers" section (previously starting at row 13) to start at row 12,
# matching the target layout.
$ws.Rows("6").Delete()

# --- Step 3: Update the "Totals:" row (now row 6) ---
$ws.Range("B6").Value2 = 7
$ws.Range("C6").Value2 = 1476

# --- Step 4: Rewrite the "Good Drivers" data rows (now rows 14-29) ---
# Pre-format the "Driver Vintage" column as text so date-like strings
# ("2024-11-10") are stored as literal text, not auto-converted to dates.
$ws.Range("E14:E29").NumberFormat = "@"

$ws.Range("A14").Value2 = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3"
$ws.Range("B14").Value2 = 34181
$ws.Range("D14").Value2 = 99.9
$ws.Range("E14").Value2 = 0

$ws.Range("A15").Value2 = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.4.0.1088"
$ws.Range("B15").Value2 = 86276
$ws.Range("D15").Value2 = 99.9
$ws.Range("E15").Value2 = 0

$ws.Range("A16").Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B16").Value2 = 56018
$ws.Range("D16").Value2 = 100
$ws.Range("E16").Value2 = 0

$ws.Range("A17").Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B17").Value2 = 34244
$ws.Range("D17").Value2 = 100
$ws.Range("E17").Value2 = 0

$ws.Range("A18").Value2 = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.3.0.1030"
$ws.Range("B18").Value2 = 17891
$ws.Range("D18").Value2 = 100
$ws.Range("E18").Value2 = 0

$ws.Range("A19").Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B19").Value2 = 442178
$ws.Range("D19").Value2 = 99.9
$ws.Range("E19").Value2 = "2024-11-10"

$ws.Range("A20").Value2 = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.0.1.5"
$ws.Range("B20").Value2 = 156943
$ws.Range("D20").Value2 = 100
$ws.Range("E20").Value2 = "2024-08-13"

$ws.Range("A21").Value2 = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5"
$ws.Range("B21").Value2 = 18738
$ws.Range("D21").Value2 = 99.9
$ws.Range("E21").Value2 = "2024-01-20"

$ws.Range("A22").Value2 = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1"
$ws.Range("B22").Value2 = 13533
$ws.Range("D22").Value2 = 100
$ws.Range("E22").Value2 = "2023-12-19"

$ws.Range("A23").Value2 = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1"
$ws.Range("B23").Value2 = 19083
$ws.Range("D23").Value2 = 100
$ws.Range("E23").Value2 = "2022-11-22"

$ws.Range("A24").Value2 = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3"
$ws.Range("B24").Value2 = 12988
$ws.Range("D24").Value2 = 100
$ws.Range("E24").Value2 = "2022-05-01"

$ws.Range("A25").Value2 = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.110.1.1"
$ws.Range("B25").Value2 = 42024
$ws.Range("D25").Value2 = 100
$ws.Range("E25").Value2 = "2022-05-01"

$ws.Range("A26").Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B26").Value2 = 77849
$ws.Range("D26").Value2 = 99.9
$ws.Range("E26").Value2 = "2021-08-18"

$ws.Range("A27").Value2 = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.70.0.6"
$ws.Range("B27").Value2 = 15504
$ws.Range("D27").Value2 = 100
$ws.Range("E27").Value2 = "2021-06-28"

$ws.Range("A28").Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B28").Value2 = 59673
$ws.Range("D28").Value2 = 100
$ws.Range("E28").Value2 = "2020-08-05"

$ws.Range("A29").Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B29").Value2 = 113652
$ws.Range("D29").Value2 = 100
$ws.Range("E29").Value2 = "2019-12-14"
